$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the "Pic" column (C) for the Francesca Findabair / Eredin rows so
# each row's image filename lines up with that row's card (Name column B),
# and normalize a couple of filenames (drop the accent in "Breacc", and
# "Daisy of The Valle" instead of "Daisy of the Valley").
$ws.Range("C12").Value = "Leaders/Francesca Findabair Pureblood Elf.png"
$ws.Range("C13").Value = "Leaders/Francesca Findabair Daisy of The Valle.png"
$ws.Range("C14").Value = "Leaders/Francesca Findabair the Beautiful.png"
$ws.Range("C15").Value = "Leaders/Francesca Findabair Queen of Dol Blathanna.png"
$ws.Range("C16").Value = "Leaders/Francesca Findabair Hope of the Aen Seidhe.png"
$ws.Range("C17").Value = "Leaders/Eredin King of the Wild Hunt.png"
$ws.Range("C18").Value = "Leaders/Eredin Commander of the Red Riders.png"
$ws.Range("C19").Value = "Leaders/Eredin Destroyer of Worlds.png"
$ws.Range("C20").Value = "Leaders/Eredin Bringer of Death.png"
$ws.Range("C21").Value = "Leaders/Eredin Breacc Glas The Treacherous.png"

# Move the active selection to C14, matching where the edit was made.
$ws.Range("C14").Select()
